$wb = $excel.ActiveWorkbook

# Sheet "建物" (building): property_category column (I) currently says "land" -> should be "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"

# Sheet "汽車" (car): property_category column (H) currently says "land" -> should be "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
